$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G gets the same width (15 chars) that the author added for the
# "FICHA" header column.
$ws.Columns.Item(7).ColumnWidth = 14.15

# Copy the formatting of the last existing header cell (F1) onto the new
# header cell (G1) so it keeps the same bold/centered/bordered style, then
# set its text.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "FICHA"

$excel.CutCopyMode = 0
